$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header capitalization: "Class_name" -> "Class_Name" (column A header, row 1)
$ws.Range("A1").Value = "Class_Name"

# Update the active selection from A3 to A2
$ws.Range("A2").Select()
